# Weekly update: insert a new "Repollo" price record (Primera/Segunda) for
# Vega Monumental Concepción dated 2021-11-11 (serial 44511) ahead of the
# existing history, shifting the previously recorded rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (rows 196:197),
# pushing the existing rows 196-233 down to 198-235.
$ws.Rows("196:197").Insert()

# Row 196: "Primera" quality entry for the new date.
$ws.Cells.Item(196, 1).Value = 11
$ws.Cells.Item(196, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(196, 3).Value = "Bíobío"
$ws.Cells.Item(196, 4).Value = 44511
$ws.Cells.Item(196, 5).Value = 8
$ws.Cells.Item(196, 6).Value = 100112006
$ws.Cells.Item(196, 7).Value = "Repollo"
$ws.Cells.Item(196, 8).Value = "Crespo record"
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 10).Value = 1500
$ws.Cells.Item(196, 11).Value = 700
$ws.Cells.Item(196, 12).Value = 800
$ws.Cells.Item(196, 13).Value = 733
$ws.Cells.Item(196, 14).Value = "$/unidad"
$ws.Cells.Item(196, 15).Value = "Región Metropolitana"
$ws.Cells.Item(196, 16).Value = 733
$ws.Cells.Item(196, 17).Value = 1
$ws.Cells.Item(196, 18).Value = "Hortaliza"

# Row 197: "Segunda" quality entry for the same new date.
$ws.Cells.Item(197, 1).Value = 11
$ws.Cells.Item(197, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(197, 3).Value = "Bíobío"
$ws.Cells.Item(197, 4).Value = 44511
$ws.Cells.Item(197, 5).Value = 8
$ws.Cells.Item(197, 6).Value = 100112006
$ws.Cells.Item(197, 7).Value = "Repollo"
$ws.Cells.Item(197, 8).Value = "Crespo record"
$ws.Cells.Item(197, 9).Value = "Segunda"
$ws.Cells.Item(197, 10).Value = 500
$ws.Cells.Item(197, 11).Value = 600
$ws.Cells.Item(197, 12).Value = 600
$ws.Cells.Item(197, 13).Value = 600
$ws.Cells.Item(197, 14).Value = "$/unidad"
$ws.Cells.Item(197, 15).Value = "Región Metropolitana"
$ws.Cells.Item(197, 16).Value = 600
$ws.Cells.Item(197, 17).Value = 1
$ws.Cells.Item(197, 18).Value = "Hortaliza"
